$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9896722435951233
$ws.Range("B1").Value = 1.627594113349915
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.579061031341553
$ws.Range("E1").Value = 1.317069053649902
